$d = $word.ActiveDocument

# Find the paragraph that holds "Docente(s) Responsável(eis)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        # Insert a brand new paragraph right after this one
        $newRange = $p.Range.InsertParagraphAfter()

        # Style the newly created paragraph as a bulleted list item
        $newPara = $p.Next()
        $newPara.Style = "ListBullet"
        $newPara.Range.Text = "11079086 - Herlandí de Souza Andrade"
        break
    }
}
